$d = $word.ActiveDocument

# Safely replace the text of a [startPos, startPos+oldLen) range with $newText.
# Inserts the new text right at the END of the old range first (a position that
# never coincides with a zero-width marker such as commentRangeStart/bookmarkStart),
# then deletes the old text from the front. This avoids the engine's quirk of
# shifting commentRangeStart/bookmarkStart markers when a replacement range's
# Start exactly matches a marker's position.
function Replace-RangeText($startPos, $oldLen, $newText) {
    $insPoint = $startPos + $oldLen
    $insRng = $d.Range($insPoint, $insPoint)
    $insRng.InsertBefore($newText)
    $oldRng = $d.Range($startPos, $insPoint)
    $oldRng.Delete()
}

# Replace every occurrence of $oldText in the document body with $newText.
function Replace-All($oldText, $newText) {
    $found = $true
    while ($found) {
        $rng = $d.Content
        $found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false)
        if ($found) {
            $s = $rng.Start
            $e = $rng.End
            $len = $e - $s
            Replace-RangeText $s $len $newText
        }
    }
}

# Header language links
Replace-All "English" "Inglés"
Replace-All " / Portuguese / French / Thai / Vietnamese / Spanish" " / Portugués / Francés / Tailandés / Vietnamita / Español"

# Brief section
Replace-All "Brief" "Breve"
Replace-All "Target audience" "Público objetivo"

# Body copy
Replace-All "We can’t wait to meet you! " "¡Estamos impacientes por conocerte! "
Replace-All "Hi " "Hola "
Replace-All "We hope you’re as excited as we are for " "Esperamos que estés tan emocionado como nosotros por el "
Replace-All "In this email, we’ve linked/attached the following documents:" "En este correo electrónico, hemos enlazado/adjuntado los siguientes documentos:"
Replace-All "Your return flight tickets" "Tus billetes de avión de ida y vuelta"
Replace-All "Your accommodation booking details" "Los datos de tu reserva de alojamiento"

Replace-All "If you have any questions, please contact us via " "Si tienes alguna pregunta, entra en contacto con nosotros por "
Replace-All " or " " o "
Replace-All "If you have any questions, please contact your country manager, " "Si tienes alguna pregunta, entra en contacto con el gestor de tu país "
Replace-All ", at " ", en "

Replace-All "See you on the " "¡Nos vemos el "
Replace-All "[DD]th" "día [DD]"

# Comment text ("choose either one" -> "elija uno de los dos")
$comment = $d.Comments.Item(3)
$comment.Range.Text = "elija uno de los dos"
